$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure the word-list table --------------------------------------
# Before:
#   A1=word    B1=def
#   A2=doggo   B2=big pupper
#   A3=pupper  B3=smol doggo
#
# After (new "pos" + "definition" + "image" columns):
#   A1=word   B1=pos   C1=definition             D1=image
#   A2=doggo  B2=noun  C2=is a full-size pupper  D2=assets/images/placeholder.png
#   A3=pupper B3=noun  C3=smol doggo             D3=assets/images/placeholder.png

$ws.Range("B1").Value = "pos"
$ws.Range("C1").Value = "definition"
$ws.Range("D1").Value = "image"

$ws.Range("B2").Value = "noun"
$ws.Range("C2").Value = "is a full-size pupper"
$ws.Range("D2").Value = "assets/images/placeholder.png"

$ws.Range("B3").Value = "noun"
$ws.Range("C3").Value = "smol doggo"
$ws.Range("D3").Value = "assets/images/placeholder.png"

# --- Column widths for the new layout --------------------------------------
# (ColumnWidth is quantized to an internal pixel grid, so these inputs are
# chosen to land as close as possible to the authored widths of
# 13.28515625 / 19 / 29.7109375.)
$ws.Columns.Item(2).ColumnWidth = 12.501302083333332
$ws.Columns.Item(3).ColumnWidth = 18.16796875
$ws.Columns.Item(4).ColumnWidth = 28.834635416666664

# --- Selection moves to B7 --------------------------------------------------
$ws.Range("B7").Select()
